# 2021.01.15 class : finished registing excel dumpfile to database
#
# The "product" sheet had its color / psize / brand columns (D, E, F)
# reordered: brand (was F) moves to D, color (was D) moves to E,
# psize (was E) moves to F. product_name / price / detail / filename
# stay put. This mirrors the user dragging the "brand" column to sit
# right after "price" (before "color") -- i.e. cutting column F and
# inserting it at column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("product")

$ws.Columns("F:F").Cut() | Out-Null
$ws.Columns("D:D").Insert() | Out-Null

# The cut/insert leaves stray default-width markers behind on the columns
# it vacated (old D, old F -- now at indices D and F again after the
# shift). Clear those columns' formatting so no bogus <col> entries show
# up, then restore the header cells' (row 1) style, which ClearFormats
# also wiped, by re-pasting formats from an untouched header cell.
$ws.Columns("D:D").ClearFormats() | Out-Null
$ws.Columns("F:F").ClearFormats() | Out-Null

$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Selection ended up on columns G:H after the reorder.
$ws.Range("G1:H1048576").Select() | Out-Null

# Tiny page-margin nudge left behind by the edit (product sheet).
$ws.PageSetup.LeftMargin = 50.380000591278076
$ws.PageSetup.RightMargin = 50.380000591278076

# Same tiny page-margin nudge on the "color" sheet (no data changes there).
$ws2 = $wb.Worksheets.Item("color")
$ws2.PageSetup.LeftMargin = 50.380000591278076
$ws2.PageSetup.RightMargin = 50.380000591278076
